$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = $null
$ws.Range("C2").Value = 26.188281108599256
$ws.Range("D2").Value = $null
$ws.Range("E2").Value = $null

$ws.Range("B3").Value = 23.27312537304482
$ws.Range("C3").Value = 29.54553601626435
$ws.Range("D3").Value = 29.129066284357918
$ws.Range("E3").Value = 14.82234274027013

$ws.Range("B1:E3").Select()
